$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in place
$ws.Range("B2").Value = 5
$ws.Range("C3").Value = 9
$ws.Range("B4").Value = 0.8
$ws.Range("C4").Value = 1.4

# Row 5 ("theta_threshold_range") is removed entirely; row 6 ("pie_threshold_range") shifts up to row 5
$ws.Rows("5").Delete()

# Update the (now shifted-up) pie_threshold_range row value
$ws.Range("C5").Value = 20

# The former row 6's B cell had a distinct font (Times New Roman); after the
# edit it matches the plain style used by the rest of the data cells.
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# Selection change as captured in the sheetView
$ws.Range("B2:C4").Select()
